$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '65.437.86'
$cell.ClearFormats()
$ws.Range('E2').Value = '  +1.26%  '

# Row 3
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '3.390.44'
$cell.ClearFormats()
$ws.Range('E3').Value = '  +0.72%  '

# Row 4
$ws.Range('E4').Value = '  -0.05%  '

# Row 5
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '559.29'
$cell.ClearFormats()
$ws.Range('E5').Value = '  -0.61%  '

# Row 6
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '175.51'
$cell.ClearFormats()
$ws.Range('E6').Value = '  -0.81%  '

# Row 7
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.632'
$cell.ClearFormats()
$ws.Range('E7').Value = '  +0.89%  '

# Row 8
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '3.378.75'
$cell.ClearFormats()
$ws.Range('E8').Value = '  +0.73%  '

# Row 9
$ws.Range('E9').Value = '  +0.02%  '

# Row 10
$ws.Range('E10').Value = '  +4.50%  '

# Row 11
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.639'
$cell.ClearFormats()
$ws.Range('E11').Value = '  +0.74%  '

# Row 12
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '53.27'
$cell.ClearFormats()
$ws.Range('E12').Value = '  -4.66%  '

# Row 13
$ws.Range('E13').Value = '  +0.41%  '

# Row 14
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '9.21'
$cell.ClearFormats()
$ws.Range('E14').Value = '  +0.83%  '

# Row 15
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '3.928.09'
$cell.ClearFormats()
$ws.Range('E15').Value = '  +0.67%  '

# Row 16
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '18.31'
$cell.ClearFormats()
$ws.Range('E16').Value = '  -0.12%  '

# Row 17
$ws.Range('E17').Value = '  +1.30%  '

# Row 18
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '3.381.69'
$cell.ClearFormats()
$ws.Range('E18').Value = '  +0.09%  '

# Row 19
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '65.389.30'
$cell.ClearFormats()
$ws.Range('E19').Value = '  +1.39%  '

# Row 20
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '11.83'
$cell.ClearFormats()
$ws.Range('E20').Value = '  -0.12%  '

# Row 21
$ws.Range('E21').Value = '  +0.91%  '

# Row 22
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '488.68'
$cell.ClearFormats()
$ws.Range('E22').Value = '  +5.80%  '

# Row 23
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '4.92'
$cell.ClearFormats()
$ws.Range('E23').Value = '  -0.54%  '

# Row 24
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '4.11'
$cell.ClearFormats()
$ws.Range('E24').Value = '  -0.69%  '

# Row 25
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '88.96'
$cell.ClearFormats()
$ws.Range('E25').Value = '  +3.01%  '

# Row 26
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '14.18'
$cell.ClearFormats()
$ws.Range('E26').Value = '  +4.39%  '

# Row 27
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '2.91'
$cell.ClearFormats()
$ws.Range('E27').Value = '  +2.14%  '

# Row 28
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '10.71'
$cell.ClearFormats()
$ws.Range('E28').Value = '  -1.39%  '

# Row 29
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '31.58'
$cell.ClearFormats()
$ws.Range('E29').Value = '  +4.43%  '

# Row 30
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '8.71'
$cell.ClearFormats()
$ws.Range('E30').Value = '  -1.46%  '

# Row 31
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '6.55'
$cell.ClearFormats()
$ws.Range('E31').Value = '  -1.97%  '

# Row 32
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '11.48'
$cell.ClearFormats()
$ws.Range('E32').Value = '  -0.42%  '

# Row 33
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '62.53'
$cell.ClearFormats()
$ws.Range('E33').Value = '  +5.52%  '

# Row 34
$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '575.48'
$cell.ClearFormats()
$ws.Range('E34').Value = '  -1.13%  '

# Row 35
$ws.Range('E35').Value = '  -0.60%  '

# Row 36
$ws.Range('E36').Value = '  +0.01%  '

# Row 37
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '3.63'
$cell.ClearFormats()
$ws.Range('E37').Value = '  +4.92%  '

# Row 38
$ws.Range('E38').Value = '  -0.03%  '

# Row 39
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '35.76'
$cell.ClearFormats()
$ws.Range('E39').Value = '  -0.73%  '

# Row 40
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.374'
$cell.ClearFormats()
$ws.Range('E40').Value = '  +0.39%  '

# Row 41
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '0.0₃0740'
$cell.ClearFormats()
$ws.Range('E41').Value = '  -2.55%  '

# Row 42
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '3.125.99'
$cell.ClearFormats()
$ws.Range('E42').Value = '  +0.98%  '

# Row 43
$ws.Range('E43').Value = '  +0.62%  '

# Row 44
$ws.Range('E44').Value = '  -1.79%  '

# Row 45
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '0.134'
$cell.ClearFormats()
$ws.Range('E45').Value = '  +0.92%  '

# Row 46
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '3.16'
$cell.ClearFormats()
$ws.Range('E46').Value = '  -1.37%  '

# Row 47
$ws.Range('E47').Value = '  -4.05%  '

# Row 48
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.ClearFormats()
$ws.Range('E48').Value = '  +0.05%  '

# Row 49
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '140.15'
$cell.ClearFormats()
$ws.Range('E49').Value = '  +2.19%  '

# Row 50
$ws.Range('E50').Value = '  -1.11%  '

# Row 51
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '8.40'
$cell.ClearFormats()
$ws.Range('E51').Value = '  -0.43%  '
